$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.252.73"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.447.80"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "2.442.46"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D17").Value = "62.130.02"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "2.443.68"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "591.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("D28").Value = "0.0₃0971"
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("D29").Value = "2.570.29"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "0.0₆0257"
$ws.Range("E48").Value = "  +18.19%  "
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0522"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
